# Scheduled market-price refresh for the Chocobo leve-profit sheets.
# Re-pulls currentAveragePrice* / Leve* figures from the latest Universalis
# snapshot and re-derives the NQ/HQ profit columns for the affected leves.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 7: The Bleat Is On
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 14: Wand-full Tonight
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 1343.1428
$ws.Range("I40").Value = 1333.6666
$ws.Range("J40").Value = 1350.25
$ws.Range("K40").Value = 1333.6666
$ws.Range("L40").Value = 1350.25
$ws.Range("M40").Value = -1158.6666
$ws.Range("N40").Value = -1700.25
# Row 135: For Tired Minds
$ws.Range("H135").Value = 607.625
$ws.Range("I135").Value = 408.7143
$ws.Range("K135").Value = 3678.4287
$ws.Range("M135").Value = -1143.4287

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 2551092
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = 3968343
$ws.Range("K2").Value = 240
$ws.Range("L2").Value = 23810058
$ws.Range("M2").Value = -127
$ws.Range("N2").Value = -23810284
# Row 4: In Hot Water
$ws.Range("H4").Value = 3015054
$ws.Range("J4").Value = 3437.6667
$ws.Range("L4").Value = 10313.0001
$ws.Range("N4").Value = -10537.0001
# Row 6: Meat-lover's Special
$ws.Range("H6").Value = 432.2
$ws.Range("I6").Value = 124.2
$ws.Range("J6").Value = 740.2
$ws.Range("K6").Value = 372.6
$ws.Range("L6").Value = 2220.6
$ws.Range("M6").Value = -259.6
$ws.Range("N6").Value = -2446.6
# Row 7: It's Always Sunny in Vylbrand
$ws.Range("H7").Value = 316.83334
$ws.Range("I7").Value = 350
$ws.Range("J7").Value = 250.5
$ws.Range("K7").Value = 1050
$ws.Range("L7").Value = 751.5
$ws.Range("M7").Value = -938
$ws.Range("N7").Value = -975.5
# Row 9: Jack of All Plates
$ws.Range("H9").Value = 1072570
$ws.Range("J9").Value = 1072570
$ws.Range("L9").Value = 3217710
$ws.Range("N9").Value = -3218158
# Row 10: A Real Fungi
$ws.Range("H10").Value = 373.66666
$ws.Range("I10").Value = 88.40000000000001
$ws.Range("J10").Value = 1800
$ws.Range("K10").Value = 265.2
$ws.Range("L10").Value = 5400
$ws.Range("M10").Value = -126.2
$ws.Range("N10").Value = -5678
# Row 12: Butter Me Up
$ws.Range("H12").Value = 131.54546
$ws.Range("I12").Value = 41.666668
$ws.Range("J12").Value = 165.25
$ws.Range("K12").Value = 125.000004
$ws.Range("L12").Value = 495.75
$ws.Range("M12").Value = 47.999996
$ws.Range("N12").Value = -841.75
# Row 13: Fishy Revelations
$ws.Range("H13").Value = 614.7273
$ws.Range("I13").Value = 110.333336
$ws.Range("J13").Value = 1220
$ws.Range("K13").Value = 331.000008
$ws.Range("L13").Value = 3660
$ws.Range("M13").Value = -163.000008
$ws.Range("N13").Value = -3996
# Row 15: Pretty Enough to Eat
$ws.Range("H15").Value = 471.2857
$ws.Range("I15").Value = 250
$ws.Range("J15").Value = 559.8
$ws.Range("K15").Value = 750
$ws.Range("L15").Value = 1679.4
$ws.Range("M15").Value = -610
$ws.Range("N15").Value = -1959.4
# Row 16: Go Ahead and Dig In
$ws.Range("H16").Value = 350
$ws.Range("I16").Value = 350
$ws.Range("K16").Value = 1050
$ws.Range("M16").Value = -877
# Row 17: Chew the Fat
$ws.Range("H17").Value = 1175.0834
$ws.Range("I17").Value = 1175.0834
$ws.Range("K17").Value = 3525.2502
$ws.Range("M17").Value = -3356.2502
# Row 19: The Bango Zango Diet
$ws.Range("H19").Value = 2650.5
$ws.Range("J19").Value = 2650.5
$ws.Range("L19").Value = 7951.5
$ws.Range("N19").Value = -8299.5
# Row 20: Omelette's Be Friends
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 21: Shy Is the Oyster
$ws.Range("H21").Value = 2659.0715
$ws.Range("I21").Value = 414.5
$ws.Range("J21").Value = 3033.1667
$ws.Range("K21").Value = 1243.5
$ws.Range("L21").Value = 9099.500100000001
$ws.Range("M21").Value = -1070.5
$ws.Range("N21").Value = -9445.500100000001
# Row 22: A Total Nut Job
$ws.Range("H22").Value = 1840.4
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 2267.3333
$ws.Range("K22").Value = 3600
$ws.Range("L22").Value = 6801.999899999999
$ws.Range("M22").Value = -3431
$ws.Range("N22").Value = -7139.999899999999
# Row 25: Flakes for Friends
$ws.Range("H25").Value = 3934.4
$ws.Range("I25").Value = 4000
$ws.Range("J25").Value = 3918
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 11754
$ws.Range("M25").Value = -11831
$ws.Range("N25").Value = -12092
# Row 26: A Grape Idea
$ws.Range("H26").Value = 12331.546
$ws.Range("I26").Value = 30112
$ws.Range("J26").Value = 2171.2856
$ws.Range("K26").Value = 90336
$ws.Range("L26").Value = 6513.8568
$ws.Range("M26").Value = -90048
$ws.Range("N26").Value = -7089.8568
# Row 27: Brain Food
$ws.Range("H27").Value = 1840.4
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 2267.3333
$ws.Range("K27").Value = 3600
$ws.Range("L27").Value = 6801.999899999999
$ws.Range("M27").Value = -3498
$ws.Range("N27").Value = -7005.999899999999
# Row 30: Picnic Panic
$ws.Range("H30").Value = 3934.4
$ws.Range("I30").Value = 4000
$ws.Range("J30").Value = 3918
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 11754
$ws.Range("M30").Value = -11898
$ws.Range("N30").Value = -11958
# Row 32: Convalescence Precedes Essence
$ws.Range("H32").Value = 3750
$ws.Range("J32").Value = 4500
$ws.Range("L32").Value = 13500
$ws.Range("N32").Value = -14066
# Row 33: Cooking with Gas
$ws.Range("H33").Value = 666784
$ws.Range("I33").Value = 2000158.2
$ws.Range("J33").Value = 96.90000000000001
$ws.Range("K33").Value = 12000949.2
$ws.Range("L33").Value = 581.4000000000001
$ws.Range("M33").Value = -12000666.2
$ws.Range("N33").Value = -1147.4
# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 127.6
$ws.Range("I38").Value = 61
$ws.Range("J38").Value = 172
$ws.Range("K38").Value = 183
$ws.Range("L38").Value = 516
$ws.Range("M38").Value = 164
$ws.Range("N38").Value = -1210
# Row 44: No More Dumpster Diving
$ws.Range("H44").Value = 2100
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 2100
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 6300
$ws.Range("N44").Value = -7096
$ws.Range("M44").ClearContents()
# Row 46: Feeding Frenzy
$ws.Range("H46").Value = 1174.2
$ws.Range("J46").Value = 1342
$ws.Range("L46").Value = 4026
$ws.Range("N46").Value = -4208
# Row 51: The Perks of Life at Sea
$ws.Range("H51").Value = 2266.6667
$ws.Range("I51").Value = 2266.6667
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 6800.000100000001
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -6340.000100000001
$ws.Range("N51").ClearContents()
# Row 57: The Egg Files
$ws.Range("H57").Value = 1233.3334
$ws.Range("I57").Value = 405
$ws.Range("K57").Value = 1215
$ws.Range("M57").Value = -656
# Row 58: Bread in the Clouds
$ws.Range("H58").Value = 1676.25
$ws.Range("I58").Value = 902.5
$ws.Range("J58").Value = 2450
$ws.Range("K58").Value = 2707.5
$ws.Range("L58").Value = 7350
$ws.Range("M58").Value = -2579.5
$ws.Range("N58").Value = -7606
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 10000864
$ws.Range("J131").Value = 929.1429000000001
$ws.Range("L131").Value = 2787.4287
$ws.Range("N131").Value = -12867.4287
# Row 132: More Mezcal
$ws.Range("H132").Value = 1926.9412
$ws.Range("J132").Value = 2382.9167
$ws.Range("L132").Value = 21446.2503
$ws.Range("N132").Value = -26506.2503

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather
$ws.Range("H122").Value = 4128.893
$ws.Range("I122").Value = 2525
$ws.Range("J122").Value = 6267.4165
$ws.Range("K122").Value = 7575
$ws.Range("L122").Value = 18802.2495
$ws.Range("M122").Value = -5125
$ws.Range("N122").Value = -23702.2495
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 2038.1837
$ws.Range("I136").Value = 1165.2106
$ws.Range("J136").Value = 5053.909
$ws.Range("K136").Value = 3495.6318
$ws.Range("L136").Value = 15161.727
$ws.Range("M136").Value = -945.6318000000001
$ws.Range("N136").Value = -20261.727
